$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProdData")
$ws.Range("A88").Value = "C28329_ValidateUserCanAddPromotionalCodeToTheConfirmAccount"
$ws.Range("B88").Value = "Refinance your RV/Camper Loan with TDECU"
$ws.Range("G88").Value = 111111
$ws.Range("H88").Value = 222222
$ws.Range("I88").Value = "vin3456789"
$ws.Range("Q88").Value = "Photo Non-Drivers License"
$ws.Range("R88").Value = "12345"
$ws.Range("S88").Value = "01012010"
$ws.Range("T88").Value = "01012025"
$ws.Range("U88").Value = "1700"
$ws.Range("V88").Value = "10"
$ws.Range("W88").Value = "11"
$ws.Range("X88").Value = "St. Mary's Hospital"
$ws.Range("Y88").Value = "2000"
$ws.Range("Z88").Value = "12"
$ws.Range("AB88").Value = "8244893"
$ws.Range("AC88").Value = "046202574"
$ws.Range("AD88").Value = "02141994"

$ws.Range("A57:AD57").Copy()
$ws.Range("A88:AD88").PasteSpecial(-4122)
